# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the Mandarina data block
# (rows 97-98), pushing the existing rows 97-133 down to 99-135.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 97:133 down by two rows to make room for the new entries.
$ws.Rows("97:98").Insert()

# --- New row 97: Murcott / Primera -----------------------------------
$ws.Cells.Item(97, 1).Value = 7
$ws.Cells.Item(97, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(97, 3).Value = "Ñuble"
$ws.Cells.Item(97, 4).Value = 44466
$ws.Cells.Item(97, 5).Value = 16
$ws.Cells.Item(97, 6).Value = "Fruta"
$ws.Cells.Item(97, 7).Value = 100102
$ws.Cells.Item(97, 8).Value = "Cítricos"
$ws.Cells.Item(97, 9).Value = 100102004
$ws.Cells.Item(97, 10).Value = "Mandarina"
$ws.Cells.Item(97, 11).Value = "Murcott"
$ws.Cells.Item(97, 12).Value = "Primera"
$ws.Cells.Item(97, 13).Value = 240
$ws.Cells.Item(97, 14).Value = 6000
$ws.Cells.Item(97, 15).Value = 6500
$ws.Cells.Item(97, 16).Value = 6250
$ws.Cells.Item(97, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(97, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(97, 19).Value = 625
$ws.Cells.Item(97, 20).Value = 10

# --- New row 98: Murcott / Segunda -------------------------------------
$ws.Cells.Item(98, 1).Value = 7
$ws.Cells.Item(98, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(98, 3).Value = "Ñuble"
$ws.Cells.Item(98, 4).Value = 44466
$ws.Cells.Item(98, 5).Value = 16
$ws.Cells.Item(98, 6).Value = "Fruta"
$ws.Cells.Item(98, 7).Value = 100102
$ws.Cells.Item(98, 8).Value = "Cítricos"
$ws.Cells.Item(98, 9).Value = 100102004
$ws.Cells.Item(98, 10).Value = "Mandarina"
$ws.Cells.Item(98, 11).Value = "Murcott"
$ws.Cells.Item(98, 12).Value = "Segunda"
$ws.Cells.Item(98, 13).Value = 120
$ws.Cells.Item(98, 14).Value = 5500
$ws.Cells.Item(98, 15).Value = 5500
$ws.Cells.Item(98, 16).Value = 5500
$ws.Cells.Item(98, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(98, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(98, 19).Value = 550
$ws.Cells.Item(98, 20).Value = 10
